$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: title change
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 32: title + link change
$ws.Range("D32").Value = "[크롤링]  What is the differences between requests and selenium?"
$ws.Range("E32").Value = "https://dodonam.tistory.com/371"

# Row 36: title + link change
$ws.Range("D36").Value = "Handling Imbalanced Data"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/366"

# Row 51: title + link change
$ws.Range("D51").Value = "[워드] 마지막 빈 페이지 삭제하기 (구역 나누기 되어 있을 경우)"
$ws.Range("E51").Value = "https://bskyvision.com/1279"
